$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Psap"
$ws.Range("C2").Value = "Gpr37l1"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 107.970093
$ws.Range("H2").Value = 323.910279
$ws.Range("I2").Value = 0.02082010292543709
$ws.Range("J2").Value = 0.02082010292543709
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4828723333333333
$ws.Range("N2").Value = 1.448617
$ws.Range("O2").Value = 0.8782126285081715
$ws.Range("P2").Value = 0.8782126285081715
$ws.Range("Q2").Value = 52.135770737127
$ws.Range("R2").Value = 469.221936634143
$ws.Range("S2").Value = 0.01828447731595878
$ws.Range("T2").Value = 0.01828447731595878

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Psap"
$ws.Range("C3").Value = "Gpr37l1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 107.970093
$ws.Range("H3").Value = 323.910279
$ws.Range("I3").Value = 0.02082010292543709
$ws.Range("J3").Value = 0.02082010292543709
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.06696300000000001
$ws.Range("N3").Value = 0.200889
$ws.Range("O3").Value = 0.1217873714918285
$ws.Range("P3").Value = 0.1217873714918285
$ws.Range("Q3").Value = 7.230001337559001
$ws.Range("R3").Value = 65.07001203803101
$ws.Range("S3").Value = 0.002535625609478312
$ws.Range("T3").Value = 0.002535625609478312

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Psap"
$ws.Range("C4").Value = "Gpr37l1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 79.28364800000001
$ws.Range("H4").Value = 237.850944
$ws.Range("I4").Value = 0.01528843465629065
$ws.Range("J4").Value = 0.01528843465629065
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4828723333333333
$ws.Range("N4").Value = 1.448617
$ws.Range("O4").Value = 0.8782126285081715
$ws.Range("P4").Value = 0.8782126285081715
$ws.Range("Q4").Value = 38.28388010493867
$ws.Range("R4").Value = 344.554920944448
$ws.Range("S4").Value = 0.01342649638527644
$ws.Range("T4").Value = 0.01342649638527644

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Psap"
$ws.Range("C5").Value = "Gpr37l1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 79.28364800000001
$ws.Range("H5").Value = 237.850944
$ws.Range("I5").Value = 0.01528843465629065
$ws.Range("J5").Value = 0.01528843465629065
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06696300000000001
$ws.Range("N5").Value = 0.200889
$ws.Range("O5").Value = 0.1217873714918285
$ws.Range("P5").Value = 0.1217873714918285
$ws.Range("Q5").Value = 5.309070921024001
$ws.Range("R5").Value = 47.781638289216
$ws.Range("S5").Value = 0.001861938271014215
$ws.Range("T5").Value = 0.001861938271014214

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Psap"
$ws.Range("C6").Value = "Gpr37l1"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1571.454671
$ws.Range("H6").Value = 4714.364013
$ws.Range("I6").Value = 0.3030269501840559
$ws.Range("J6").Value = 0.3030269501840559
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4828723333333333
$ws.Range("N6").Value = 1.448617
$ws.Range("O6").Value = 0.8782126285081715
$ws.Range("P6").Value = 0.8782126285081715
$ws.Range("Q6").Value = 758.8119837133357
$ws.Range("R6").Value = 6829.307853420021
$ws.Range("S6").Value = 0.2661220944299545
$ws.Range("T6").Value = 0.2661220944299545

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Psap"
$ws.Range("C7").Value = "Gpr37l1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1571.454671
$ws.Range("H7").Value = 4714.364013
$ws.Range("I7").Value = 0.3030269501840559
$ws.Range("J7").Value = 0.3030269501840559
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.06696300000000001
$ws.Range("N7").Value = 0.200889
$ws.Range("O7").Value = 0.1217873714918285
$ws.Range("P7").Value = 0.1217873714918285
$ws.Range("Q7").Value = 105.229319134173
$ws.Range("R7").Value = 947.0638722075571
$ws.Range("S7").Value = 0.03690485575410141
$ws.Range("T7").Value = 0.03690485575410141

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Psap"
$ws.Range("C8").Value = "Gpr37l1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 35.82847833333333
$ws.Range("H8").Value = 107.485435
$ws.Range("I8").Value = 0.006908881763784279
$ws.Range("J8").Value = 0.00690888176378428
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4828723333333333
$ws.Range("N8").Value = 1.448617
$ws.Range("O8").Value = 0.8782126285081715
$ws.Range("P8").Value = 0.8782126285081715
$ws.Range("Q8").Value = 17.30058093259944
$ws.Range("R8").Value = 155.705228393395
$ws.Range("S8").Value = 0.006067467213825164
$ws.Range("T8").Value = 0.006067467213825165

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Psap"
$ws.Range("C9").Value = "Gpr37l1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 35.82847833333333
$ws.Range("H9").Value = 107.485435
$ws.Range("I9").Value = 0.006908881763784279
$ws.Range("J9").Value = 0.00690888176378428
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.06696300000000001
$ws.Range("N9").Value = 0.200889
$ws.Range("O9").Value = 0.1217873714918285
$ws.Range("P9").Value = 0.1217873714918285
$ws.Range("Q9").Value = 2.399182394635
$ws.Range("R9").Value = 21.592641551715
$ws.Range("S9").Value = 0.0008414145499591151
$ws.Range("T9").Value = 0.0008414145499591151

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "Psap"
$ws.Range("C10").Value = "Gpr37l1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 929.2825723333332
$ws.Range("H10").Value = 2787.847717
$ws.Range("I10").Value = 0.1791955370715012
$ws.Range("J10").Value = 0.1791955370715012
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4828723333333333
$ws.Range("N10").Value = 1.448617
$ws.Range("O10").Value = 0.8782126285081715
$ws.Range("P10").Value = 0.8782126285081715
$ws.Range("Q10").Value = 448.7248440285987
$ws.Range("R10").Value = 4038.523596257389
$ws.Range("S10").Value = 0.1573717836284966
$ws.Range("T10").Value = 0.1573717836284966

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Psap"
$ws.Range("C11").Value = "Gpr37l1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 929.2825723333332
$ws.Range("H11").Value = 2787.847717
$ws.Range("I11").Value = 0.1791955370715012
$ws.Range("J11").Value = 0.1791955370715012
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.06696300000000001
$ws.Range("N11").Value = 0.200889
$ws.Range("O11").Value = 0.1217873714918285
$ws.Range("P11").Value = 0.1217873714918285
$ws.Range("Q11").Value = 62.227548891157
$ws.Range("R11").Value = 560.047940020413
$ws.Range("S11").Value = 0.02182375344300464
$ws.Range("T11").Value = 0.02182375344300464

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Psap"
$ws.Range("C12").Value = "Gpr37l1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2462.03833
$ws.Range("H12").Value = 7386.11499
$ws.Range("I12").Value = 0.4747600933989308
$ws.Range("J12").Value = 0.4747600933989309
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.4828723333333333
$ws.Range("N12").Value = 1.448617
$ws.Range("O12").Value = 0.8782126285081715
$ws.Range("P12").Value = 0.8782126285081715
$ws.Range("Q12").Value = 1188.850193163203
$ws.Range("R12").Value = 10699.65173846883
$ws.Range("S12").Value = 0.41694030953466
$ws.Range("T12").Value = 0.4169403095346601

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Psap"
$ws.Range("C13").Value = "Gpr37l1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2462.03833
$ws.Range("H13").Value = 7386.11499
$ws.Range("I13").Value = 0.4747600933989308
$ws.Range("J13").Value = 0.4747600933989309
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.06696300000000001
$ws.Range("N13").Value = 0.200889
$ws.Range("O13").Value = 0.1217873714918285
$ws.Range("P13").Value = 0.1217873714918285
$ws.Range("Q13").Value = 164.86547269179
$ws.Range("R13").Value = 1483.78925422611
$ws.Range("S13").Value = 0.05781978386427077
$ws.Range("T13").Value = 0.05781978386427077
